# Reformat the 'Programa' and 'Bibliografia' body paragraphs: both were
# authored as a single run containing every numbered/bulleted item back
# to back; split each item onto its own line with a manual line break
# (<w:br/>) between <w:t> runs, matching the canonical OOXML.
$d = $word.ActiveDocument

function Set-ParagraphLines {
    param(
        [int]$ParaIndex,
        [string]$ExpectedStart,
        [string]$ParagraphXml
    )
    $range = $d.Paragraphs($ParaIndex).Range
    if (-not $range.Text.StartsWith($ExpectedStart)) {
        throw "Paragraph $ParaIndex does not start with the expected text; aborting."
    }
    # Replace the whole paragraph (InsertXML on a Range spanning it swaps
    # its content in place) with the pre-built run/break XML below so the
    # resulting run/break structure -- including xml:space=preserve on runs
    # with significant leading/trailing whitespace -- matches exactly,
    # rather than relying on Find/Replace's own serialization.
    [void]$range.InsertXML($ParagraphXml)
}

$programaXml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t>1. Sistemas de Informação</w:t><w:br/><w:t>1.1. Sistemas de Processamento de Informações;</w:t><w:br/><w:t>1.2. Sistemas de Informações Gerenciais;</w:t><w:br/><w:t>1.3. Sistema de Apoio à Decisão;</w:t><w:br/><w:t>1.4. Sistemas de Informação no Comércio Eletrônico;</w:t><w:br/><w:t>1.5. Sistemas de Informação em Cadeia de Suprimentos;</w:t><w:br/><w:t>1.6. Sistemas inteligentes nos negócios;</w:t><w:br/><w:t xml:space=`"preserve`">1.7. Sistemas estratégicos. </w:t><w:br/><w:t>2. Projeto de Sistemas de Informação.</w:t><w:br/><w:t>2.1. Especificação das Saídas;</w:t><w:br/><w:t>2.2. Especificação dos Arquivos;</w:t><w:br/><w:t>2.3. Especificação das Entradas;</w:t><w:br/><w:t>2.4. Especificação do Processamento.</w:t><w:br/><w:t>3. Tecnologia de Informação.</w:t><w:br/><w:t>3.1. Evolução da Computação;</w:t><w:br/><w:t>3.2. Recursos Computacionais.</w:t><w:br/><w:t>4. Processo de Desenvolvimento de Sistemas de Informação.</w:t><w:br/><w:t>4.1. Definição do Negócio;</w:t><w:br/><w:t>4.2. Identificação do Problema e/ou Oportunidades;</w:t><w:br/><w:t>4.3. Seleção do Sistema de Informação;</w:t><w:br/><w:t>4.4. Implementação do Sistema de Informação;</w:t><w:br/><w:t>4.5. Avaliação da Eficácia do Sistema de Informação;</w:t></w:r></w:p>"
Set-ParagraphLines 12 "1. Sistemas de Informação" $programaXml

$biblioXml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t>HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004</w:t><w:br/><w:t>LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gestão integrada de processos e da tecnologia da informação. São Paulo:Atlas, 2006.</w:t><w:br/><w:t>LAURINDO, F.J.B. Tecnologia da Informação: Eficácia nas Organizações. São Paulo, Editora Futura, 2002.</w:t><w:br/><w:t>STAIR, R.M., Princípios de Sistema de Informação: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.</w:t><w:br/><w:t>TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.</w:t><w:br/><w:t>TURBAN, E., RAIANER JR, K., POTTER, R. E., Administração de Tecnologia da Informação: Teoria e Prática”, São Paulo, Editora Campus, 2003.</w:t></w:r></w:p>"
Set-ParagraphLines 16 "HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004" $biblioXml


Write-Output "Programa paragraph now has $($d.Paragraphs(12).Range.Text.Split([char]11).Length) lines"
Write-Output "Bibliografia paragraph now has $($d.Paragraphs(16).Range.Text.Split([char]11).Length) lines"
